$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 (A12) currently carries the "date only" number format because it
# used to be the last row of the series. Capture that format so it can be
# moved onto the new last row (row 13).
$lastRowDateFormat = $ws.Range("A12").NumberFormat

# Every other date cell (e.g. A2) uses the regular "date + time" number
# format; row 12 should switch to that now that it is no longer last.
$ws.Range("A12").NumberFormat = $ws.Range("A2").NumberFormat

# Append the new day's data.
$ws.Range("A13").Value = 45753
$ws.Range("B13").Value = 46
$ws.Range("C13").Value = 51
$ws.Range("D13").Value = 46

# New last row gets the "date only" format previously used by A12.
$ws.Range("A13").NumberFormat = $lastRowDateFormat
